$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New part-number row (was previously occupied by the "Total" label):
# Force column A to be stored as text so "111954120" becomes a shared string
# (matching the existing text part numbers in A2/A3) instead of a number,
# then drop the temporary number-format style so no stray formatting is left
# on the cell.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "111954120"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = 0.1
$ws.Range("C4").Value = 509.97

# The "Total" label moves down to row 5 and is updated to reflect the
# new sum: 269.7 + 4794 + 509.97 = 5573.67
$ws.Range("C5").Value = "Total: 5573.67"
